
# "Update on DB connect" — append new REP/Bee-No rows to the Accounts sheet
# and resize its columns to fit the new, wider data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Accounts")

# --- Resize the columns (A:40, B:6, C:12, D:32 "characters") ---------------
# Excel's ColumnWidth property is expressed in "characters of the Normal
# style font" and round-trips to the stored <col width="..."> value with a
# constant +5/6 offset (the internal pixel-padding term). Subtracting that
# offset here makes the value actually written to the sheet come out exact.
$padding = 5 / 6
$ws.Columns("A:A").ColumnWidth = 40 - $padding
$ws.Columns("B:B").ColumnWidth = 6 - $padding
$ws.Columns("C:C").ColumnWidth = 12 - $padding
$ws.Columns("D:D").ColumnWidth = 32 - $padding

# --- New rows of rep/account data ------------------------------------------
$rows = @(
    @("TylerLarson@rep.com",     "REP", "400463", "No Cost Signup"),
    @("StephanieLynn@rep.com",   "REP", "400465", "No Cost Signup"),
    @("TinaJohnson@rep.com",     "REP", "400466", "No Cost Signup"),
    @("MaryJoyce@rep.com",       "REP", "400470", "30$ kit Signup"),
    @("MichaelFisher@rep.com",   "REP", "400472", "30$ kit Signup"),
    @("GeorgeKeller@rep.com",    "REP", "400478", "30$ kit Signup"),
    @("StacyJenkins@rep.com",    "REP", "400479", "No Cost Signup")
)

$r = 5
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row[0]

    $ws.Range("B$r").Value = $row[1]

    # Force the Bee No. to be stored as text, not a number, matching the
    # existing C2:C4 cells.
    $ws.Range("C$r").NumberFormat = "@"
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("C$r").Style = "Normal"

    # The Notes column values are literally wrapped in single quotes
    # (e.g. 'No Cost Signup'). A leading `'` fed through .Value is treated
    # by Excel as a quote-prefix marker rather than literal text, so a
    # second leading `'` is used to get one literal leading quote through.
    $ws.Range("D$r").Value = "''" + $row[3] + "'"
    $ws.Range("D$r").Style = "Normal"

    $r++
}
